$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill A1:D1 with the text "yo" (becomes a shared string entry)
$ws.Range("A1:D1").Value = "yo"

# Match the author's final cursor position/selection (F13)
$ws.Range("F13").Select()
